$d = $word.ActiveDocument

# 1. Rename the TOC bookmark on the chapter title (id stays 0, name changes).
#    Bookmark.Name is not writable in this host, so delete + re-add at the
#    same (zero-width) range to preserve its position.
$bm = $d.Bookmarks("_Toc442187695")
$bmRange = $bm.Range
$bmStart = $bmRange.Start
$bmEnd = $bmRange.End
$bm.Delete()
$d.Bookmarks.Add("_Toc445473431", $d.Range($bmStart, $bmEnd))

# 2. "...radiance from 650 to 950 nm in wavelength in two dimensional..."
#    -> "...radiance from 650 to 950 nm in two dimensional..."
$d.Content.Find.Execute("in wavelength in two dimensional", $true, $false, $false, $false, $false, `
    $true, 1, $false, "in two dimensional", 2) | Out-Null

# 3. "...nominal 210 m both vertical and horizontal direction."
#    -> "...nominal 210 m both in the vertical and horizontal directions."
$d.Content.Find.Execute("both vertical and horizontal direction.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "both in the vertical and horizontal directions.", 2) | Out-Null

# 4. "...can be increased to compensate the decrease in radiance."
#    -> "...can be increased to compensate for the overall decrease in radiance."
$d.Content.Find.Execute("compensate the decrease in radiance", $true, $false, $false, $false, $false, `
    $true, 1, $false, "compensate for the overall decrease in radiance", 2) | Out-Null

# 5. "...agreed well to the nearest OSIRIS scans but had some disagreement in extinction values..."
#    -> "...but had some large discrepancies in extinction values..."
$d.Content.Find.Execute("some disagreement in extinction values", $true, $false, $false, $false, $false, `
    $true, 1, $false, "some large discrepancies in extinction values", 2) | Out-Null

# 6. "...calibrate DC offset and dark current change during the flight..."
#    -> "...calibrate DC offset and dark current changes during the flight..."
$d.Content.Find.Execute("dark current change during the flight", $true, $false, $false, $false, $false, `
    $true, 1, $false, "dark current changes during the flight", 2) | Out-Null
